# factor(tests): write classes to factor first two tests
#
# This script:
#   1. Removes the extra "G" column of data from the "expected" sheet
#      (it becomes the new "expected2" sheet's extra column instead).
#   2. Adds a new worksheet named "expected2" at the end of the workbook,
#      containing the same layout as "Feuille2" (columns A-E) plus the
#      extra column (F) that used to live in column G of "expected".
#   3. Restores per-sheet selections / active-cell bookkeeping and leaves
#      "expected2" as the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuille2")
$ws2 = $wb.Worksheets.Item("expected")

# --- 1. Strip column G from "expected"; the data moves to the new sheet ---
[void]$ws2.Range("G1:G9").Delete()

# --- 2. Create "expected2" as the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "expected2"

# Header row
$ws3.Range("A1").Value = "Evaluatisme"
$ws3.Range("B1").Value = "Absolutisme"
$ws3.Range("C1").Value = "Multiplisme"
$ws3.Range("E1").Value = "Colonne de(s) maximum(s)"

# Row 2
$ws3.Range("A2").NumberFormat = "General"
$ws3.Range("B2").Value = "prime"
$ws3.Range("C2").Value = "congruent1_o.jpg"
$ws3.Range("D2").Value = "prime_congruent"
$ws3.Range("E2").Value = "Congruent1;o;.jpg"
$ws3.Range("F2").Value = "o"

# Row 3
$ws3.Range("A3").NumberFormat = "General"
$ws3.Range("B3").Value = "probe"
$ws3.Range("C3").Value = "incongruent1_n.jpg"
$ws3.Range("D3").Value = "probe_incongruent"
$ws3.Range("E3").Value = "Incongruent1;n.jpg"
$ws3.Range("F3").Value = "n.jpg"

# Row 4
$ws3.Range("A4").NumberFormat = "General"
$ws3.Range("B4").Value = "qfqfq"
$ws3.Range("C4").Value = "congruent2_o.jpg"
$ws3.Range("E4").Value = "Congruent2;o.jpg"
$ws3.Range("F4").Value = "o.jpg"

# Row 5
$ws3.Range("A5").NumberFormat = "General"
$ws3.Range("B5").Value = "qvfqf"
$ws3.Range("C5").Value = "congruent3_n.jpg"
$ws3.Range("E5").Value = "congruent3;n.jpg"
$ws3.Range("F5").Value = "n.jpg"

# Row 6
$ws3.Range("B6").Value = "probe"
$ws3.Range("C6").Value = "congruent5_o.jpg"
$ws3.Range("D6").Value = "probe_congruent"
$ws3.Range("E6").Value = "Congruent5;o.jpg"
$ws3.Range("F6").Value = "o.jpg"

# Row 7
$ws3.Range("B7").Value = "pqfqf"
$ws3.Range("C7").Value = "d1_o.jpg"
$ws3.Range("E7").Value = "d1_o;.jpg"
$ws3.Range("F7").Value = ".jpg"

# Row 8
$ws3.Range("B8").Value = "Prime"
$ws3.Range("C8").Value = "congruent1_o.jpg"
$ws3.Range("E8").Value = "congruent1;_o.jpg"
$ws3.Range("F8").Value = "_o.jpg"

# Row 9
$ws3.Range("B9").Value = "prime"
$ws3.Range("C9").Value = "neutre12_n.jpg"
$ws3.Range("D9").Value = "prime_neutre"
$ws3.Range("E9").Value = "Neutre1;2_n.jpg"
$ws3.Range("F9").Value = "2_n.jpg"

# --- 3. Selections / active sheet bookkeeping ---
[void]$ws1.Range("D2").Select()
[void]$ws2.Range("G1").Select()
[void]$ws3.Range("D1").Select()
[void]$ws3.Activate()
